$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy header style from an existing header cell (e.g., E1) to the new headers
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Boolean data values for rows 2-8
$values = @(
    @($false, $false, $false),
    @($false, $true,  $true),
    @($false, $true,  $false),
    @($false, $false, $false),
    @($false, $false, $true),
    @($false, $false, $false),
    @($false, $false, $false)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 6).Value = $values[$i][0]
    $ws.Cells.Item($row, 7).Value = $values[$i][1]
    $ws.Cells.Item($row, 8).Value = $values[$i][2]
}
